# productListingxlsx.xlsx - "excel and csv update"
#
# The product-image URLs for the two product rows that pointed at
# "Red Jacket.jpeg" / "White Shirt.jpeg" are repointed to the new GitHub
# hosted locations for the corresponding images. All other data is
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "https://github.com/Vincent-Cayadi/Website-Backup/blob/main/productImages/Black%20Sweater.jpeg"
$ws.Range("D3").Value = "https://github.com/Vincent-Cayadi/Website-Backup/blob/main/productImages/White%20Shirt.jpeg"

# Leave the selection where the author last left it before saving.
$ws.Range("D3").Select() | Out-Null
